# Update crypto price/volume table cells per the Nov 16 2023 GitHub Actions refresh.
# Values that look numeric are written with a leading apostrophe so Excel stores
# them as text (quotePrefix), matching the source data's text-formatted Price column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.038.95"
$ws.Range("E2").Value = "  -4.15%  "

$ws.Range("D3").Value = "1.956.17"
$ws.Range("E3").Value = "  -4.06%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'241.13"
$ws.Range("E5").Value = "  -4.00%  "

$ws.Range("D6").Value = "'0.625"
$ws.Range("E6").Value = "  -3.38%  "

$ws.Range("D7").Value = "'60.18"
$ws.Range("E7").Value = "  -9.28%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "'0.370"
$ws.Range("E9").Value = "  -1.26%  "

$ws.Range("D10").Value = "'56.31"
$ws.Range("E10").Value = "  -5.36%  "

$ws.Range("D11").Value = "'0.0800"
$ws.Range("E11").Value = "  +6.23%  "

$ws.Range("D12").Value = "'0.103"
$ws.Range("E12").Value = "  -0.97%  "

$ws.Range("D13").Value = "'0.847"
$ws.Range("E13").Value = "  -6.25%  "

$ws.Range("D14").Value = "'13.91"
$ws.Range("E14").Value = "  -7.80%  "

$ws.Range("D15").Value = "'21.67"
$ws.Range("E15").Value = "  +4.72%  "

$ws.Range("D16").Value = "2.239.83"
$ws.Range("E16").Value = "  -4.25%  "

$ws.Range("D17").Value = "'5.38"
$ws.Range("E17").Value = "  -3.89%  "

$ws.Range("D18").Value = "1.955.73"
$ws.Range("E18").Value = "  -4.20%  "

$ws.Range("D19").Value = "35.894.27"
$ws.Range("E19").Value = "  -4.17%  "

$ws.Range("D20").Value = "'70.73"
$ws.Range("E20").Value = "  -3.48%  "

$ws.Range("D21").Value = "0.0₃0850"
$ws.Range("E21").Value = "  -2.66%  "

$ws.Range("D22").Value = "'235.07"
$ws.Range("E22").Value = "  -0.98%  "

$ws.Range("D23").Value = "'5.17"
$ws.Range("E23").Value = "  -3.46%  "

$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("E25").Value = "  -5.47%  "

$ws.Range("D26").Value = "'2.28"
$ws.Range("E26").Value = "  -4.04%  "

$ws.Range("D27").Value = "'9.70"
$ws.Range("E27").Value = "  +1.40%  "

$ws.Range("D28").Value = "'159.09"
$ws.Range("E28").Value = "  -3.68%  "

$ws.Range("D29").Value = "'19.73"
$ws.Range("E29").Value = "  -0.88%  "

$ws.Range("D30").Value = "'0.130"
$ws.Range("E30").Value = "  +17.68%  "

$ws.Range("D31").Value = "'0.119"
$ws.Range("E31").Value = "  -2.18%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.82"
$ws.Range("E32").Value = "  -7.69%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.13"
$ws.Range("E33").Value = "  -6.68%  "

$ws.Range("D34").Value = "'0.0615"
$ws.Range("E34").Value = "  +0.55%  "

$ws.Range("D35").Value = "'4.37"
$ws.Range("E35").Value = "  -7.51%  "

$ws.Range("D36").Value = "'6.25"
$ws.Range("E36").Value = "  +2.93%  "

$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").Value = "'1.82"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("E39").Value = "  -8.37%  "

$ws.Range("D40").Value = "'3.04"
$ws.Range("E40").Value = "  +11.08%  "

$ws.Range("D41").Value = "'0.0979"
$ws.Range("E41").Value = "  -4.51%  "

$ws.Range("E42").Value = "  -1.35%  "

$ws.Range("E43").Value = "  -3.52%  "

$ws.Range("D44").Value = "'0.0210"
$ws.Range("E44").Value = "  -3.54%  "

$ws.Range("E45").Value = "  -4.90%  "

$ws.Range("D46").Value = "'91.80"
$ws.Range("E46").Value = "  -3.46%  "

$ws.Range("D47").Value = "'15.94"
$ws.Range("E47").Value = "  -5.93%  "

$ws.Range("E48").Value = "  -7.87%  "

$ws.Range("D49").Value = "1.329.76"
$ws.Range("E49").Value = "  -6.77%  "

$ws.Range("E50").Value = "  -5.27%  "

$ws.Range("D51").Value = "2.134.83"
$ws.Range("E51").Value = "  -4.15%  "
